# Update row 2 ("ig_traj_460" attributions) with the new relative-direction
# values. Only the cells that actually changed vs. the original are touched;
# every other cell (rows 1, 3 and the untouched row-2 columns) is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1196116061141161
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.290569133234185
$ws.Range("E2").Value = 0.006638047860872211
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.08015518083814908
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2541322985548342
$ws.Range("N2").Value = 0.02476157060175758
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1114215457575703
$ws.Range("V2").Value = 0.01317060550227939
$ws.Range("W2").Value = -0.01740300039027242
$ws.Range("Y2").Value = -0
$ws.Range("Z2").Value = -0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.007358126984907245
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0.0209489242586958
$ws.Range("AF2").Value = -0.001278400681974891
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.05293777923492834
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.04805918661798327
$ws.Range("AO2").Value = 0.07416228179326181
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AS2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1561542071371479
$ws.Range("AW2").Value = 0.1159826283168068
$ws.Range("AX2").Value = -0.02550510037105823
$ws.Range("AY2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.01975330946610538
$ws.Range("BE2").Value = -0
$ws.Range("BF2").Value = 0.1173026497105192
$ws.Range("BG2").Value = 0.01848489159634745
$ws.Range("BI2").Value = 0
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.01234147895666798
$ws.Range("BO2").Value = -0.04048451837496309
$ws.Range("BP2").Value = -0.04534039424958565
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.05915712025865116
$ws.Range("BW2").Value = 0
$ws.Range("BX2").Value = 0.02790266473584342
$ws.Range("BY2").Value = -0.0324154409580402
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.02752952368981012
$ws.Range("CF2").Value = -0
$ws.Range("CG2").Value = -0.04935487713575477
$ws.Range("CH2").Value = 0.01352769766097447
$ws.Range("CI2").Value = 0
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.02357779391119231
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.01519330751800707
$ws.Range("CQ2").Value = 0.0674696792572889
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04430740944458727
$ws.Range("CY2").Value = -0.04969995472952137
$ws.Range("CZ2").Value = 0.004314839788907577
$ws.Range("DA2").Value = -0
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.04071710234245642
$ws.Range("DH2").Value = 0.006051752272508887
$ws.Range("DI2").Value = 0.05840709464602017
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.02989034350703162
$ws.Range("DP2").Value = -0
$ws.Range("DQ2").Value = 0.04157758598046472
$ws.Range("DR2").Value = -0.06208716646384525
$ws.Range("DS2").Value = -0
$ws.Range("DT2").Value = 0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05999670130279144
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.005262735983777017
$ws.Range("EA2").Value = -0.03917335870846283
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04390391319024155
$ws.Range("EI2").Value = 0.09826642592587737
$ws.Range("EJ2").Value = -0.01862336688472342
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.05288800328953746
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.05981517127729411
$ws.Range("ES2").Value = 0.0196735427834946
$ws.Range("ET2").Value = 0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.0387038909027299
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = -0.04123072011519976
$ws.Range("FB2").Value = 0.01666623914432457
$ws.Range("FC2").Value = 0
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.009345571191879251
$ws.Range("FJ2").Value = -0.01487059634927098
$ws.Range("FK2").Value = 0.02129248369880309
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.002131586852149674
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = 0.009157938490005671
$ws.Range("FT2").Value = -0.002676297973793981
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.0440100622146601
$ws.Range("GB2").Value = 0.01701625603450851
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
